# "viz y tablas update"
# Sheet "Ficha técnica": remove the DIMENSIÓN / Disponibilidad row, and
# append two new rows (TIPOIND / Resultados, CITA / UMAD con base en
# Observatorio Territorio Uruguay - OPP) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

# Row 3 currently holds "DIMENSIÓN" / "Disponibilidad" -- delete the whole
# row so everything below shifts up one position.
$ws.Rows.Item(3).Delete()

# Append the two new rows after the current last row (now row 6: CÁLCULO).
$ws.Cells.Item(7, 1).Value = "TIPOIND"
$ws.Cells.Item(7, 2).Value = "Resultados"

$ws.Cells.Item(8, 1).Value = "CITA"
$ws.Cells.Item(8, 2).Value = "UMAD con base en Observatorio Territorio Uruguay - OPP"
